# HW1_BoxModel_Hsieh.xlsx - "Added headings to my figures"
#
# The underlying edit re-classifies the soil-type (zone) selector for the
# top five internal nodes of the 1-D box model on the "model and key plot"
# sheet: cells G8:G11 move from zone 3 to zone 2, and G12:G16 move from
# zone 3 to zone 1. Every other number on the sheet (H, I, J, D7:D9, C11,
# C12, L/M/N flags, the chart caches, etc.) is a downstream formula result
# that recalculates automatically - including the iterative (circular
# reference) head solve the sheet relies on (Excel: iterative calc is
# already enabled on this workbook).
#
# The view also changed: the sheet is now shown zoomed to 80% with the
# selection parked on AB15 (and scrolled down a bit) instead of K10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model and key plot")

# --- 1. Re-assign the soil/K zone for nodes 8-16 -------------------------
# (G17:G20 already hold zone 3 and are left untouched.)
$newZones = @{
    8  = 2
    9  = 2
    10 = 2
    11 = 2
    12 = 1
    13 = 1
    14 = 1
    15 = 1
    16 = 1
}

# This sheet has a circular reference (node heads I9:I19 depend on their
# neighbors, which depend on them back) solved through Excel's iterative
# calculation. Each recalculation pass only advances the Gauss-Seidel
# solve a little, so re-apply the inputs and recalculate repeatedly until
# the dependent head/flux columns settle down, just like Excel converging
# over its 100-iteration default.
for ($pass = 0; $pass -lt 60; $pass++) {
    foreach ($row in $newZones.Keys) {
        $ws.Cells.Item($row, 7).Value = $newZones[$row]
    }
    $excel.Calculate()
}

# --- 2. Update the sheet view (zoom + selection) -------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 80
$ws.Range("AB15").Select()
